# Replace the old "Donations by cause" breakdown table (Academic research &
# education / Social causes / Arts / Politics / Sport / Total amount
# donated) with the new "training & follow-up courses" breakdown table, and
# drop the trailing columns (old Sport / Total amount donated) that have no
# counterpart in the new data, shrinking the used range back to A1:E4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B1:E4 hold numeric-looking / comma-grouped strings ("362", "2,082", ...)
# that must stay literal text, not get auto-converted to real numbers, so
# force text format before writing them in.
$ws.Range("B1:E4").NumberFormat = "@"

# Written column-major (all of B top-to-bottom, then all of C, ...) so the
# shared-string table is populated in the same order Excel itself would use.
$ws.Range("B1").Value = "Prevention of money laundering (incl. follow-up courses)"
$ws.Range("B2").Value = "362"
$ws.Range("B3").Value = "2,082"
$ws.Range("B4").Value = "3,206"

$ws.Range("C1").Value = "Ban on market price manipulation (incl. follow-up courses)*"
$ws.Range("C2").Value = "494"
$ws.Range("C3").Value = "351"
$ws.Range("C4").Value = "468"

$ws.Range("D1").Value = "Fraud prevention for managers (incl. follow-up courses)"
$ws.Range("D2").Value = "169"
$ws.Range("D3").Value = "447"
$ws.Range("D4").Value = "125"

$ws.Range("E1").Value = "Fraud prevention for employees (incl. follow-up courses)"
$ws.Range("E2").Value = "1,153"
$ws.Range("E3").Value = "3,920"
$ws.Range("E4").Value = "728"

# Drop the temporary text-number-format override again so the cells keep
# their original (default) style once the text is safely stored.
$ws.Range("B1:E4").ClearFormats()

# The old "Arts/Politics/Sport" + "Total amount donated (€)" columns (F, G)
# have no equivalent in the new table -- remove them entirely.
$ws.Range("F1:G4").ClearContents()
